$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Tags" column (column D) entirely - shifts everything right of it
# one column to the left (E->D, F->E, G->F, H->G, I->H).
$ws.Columns("D").Delete()

# Insert a new row for problem 584 "Find Customer Referee" at row 4, pushing
# the existing placeholder rows down by one.
$ws.Rows("4").Insert()

# The list still only has 16 total rows, so drop the now-duplicated trailing
# placeholder row that got pushed down to row 17.
$ws.Rows("17").Delete()

# Fill in the new row 4 with the finished-problem data.
$checkMark = [char]0x2705
$apos = [char]0x2019

$ws.Range("A4").Value2 = 584
$ws.Range("B4").Value2 = "Find Customer Referee"
$ws.Range("C4").Value2 = "SELECT"
$ws.Range("E4").Value2 = "Easy"
$ws.Range("F4").Value2 = 1
$ws.Range("G4").Value2 = $checkMark
$ws.Range("H4").Value2 = "Given 1 sol and didn${apos}t see solutions"

$ws.Rows("4").RowHeight = 30

# Restore the active selection shown in the author's session.
$ws.Range("H5").Select()
